$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35, shifting existing rows 35:75 down to 36:75
$ws.Rows.Item(35).Insert()

# Populate the new row 35 with the new data entry
$ws.Cells.Item(35, 1).Value = 6
$ws.Cells.Item(35, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(35, 3).Value = "Metropolitana"
$ws.Cells.Item(35, 4).Value = 44781
$ws.Cells.Item(35, 5).Value = 13
$ws.Cells.Item(35, 6).Value = "Fruta"
$ws.Cells.Item(35, 7).Value = 100108
$ws.Cells.Item(35, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(35, 9).Value = 100108007
$ws.Cells.Item(35, 10).Value = "Coco"
$ws.Cells.Item(35, 11).Value = "Sin especificar"
$ws.Cells.Item(35, 12).Value = "Primera"
$ws.Cells.Item(35, 13).Value = 200
$ws.Cells.Item(35, 14).Value = 23000
$ws.Cells.Item(35, 15).Value = 24000
$ws.Cells.Item(35, 16).Value = 23500
$ws.Cells.Item(35, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(35, 18).Value = "Perú"
$ws.Cells.Item(35, 19).Value = 1175
$ws.Cells.Item(35, 20).Value = 20
